$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 193. This pushes the existing rows
# 193, 194, 195 down to 195, 196, 197 - keeping all of their original
# content and formatting intact.
$ws.Rows("193:194").Insert()

# --- New row 193 ---
$ws.Cells.Item(193, 1).Value = 1
$ws.Cells.Item(193, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(193, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(193, 4).Value = 45239
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(193, 5).Value = 15
$ws.Cells.Item(193, 6).Value = 100112008
$ws.Cells.Item(193, 7).Value = "Coliflor"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Segunda"
$ws.Cells.Item(193, 10).Value = 450
$ws.Cells.Item(193, 11).Value = 500
$ws.Cells.Item(193, 12).Value = 600
$ws.Cells.Item(193, 13).Value = 522
$ws.Cells.Item(193, 14).Value = "$/unidad"
$ws.Cells.Item(193, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(193, 16).Value = 522
$ws.Cells.Item(193, 17).Value = 1
$ws.Cells.Item(193, 18).Value = "Hortaliza"

# --- New row 194 ---
$ws.Cells.Item(194, 1).Value = 1
$ws.Cells.Item(194, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(194, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(194, 4).Value = 45239
$ws.Cells.Item(194, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(194, 5).Value = 15
$ws.Cells.Item(194, 6).Value = 100112008
$ws.Cells.Item(194, 7).Value = "Coliflor"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Tercera"
$ws.Cells.Item(194, 10).Value = 659
$ws.Cells.Item(194, 11).Value = 300
$ws.Cells.Item(194, 12).Value = 400
$ws.Cells.Item(194, 13).Value = 354
$ws.Cells.Item(194, 14).Value = "$/unidad"
$ws.Cells.Item(194, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(194, 16).Value = 354
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = "Hortaliza"
